$d = $word.ActiveDocument

# The paragraph currently reads (as three separately-formatted runs, all
# colored C00000):
#   "(This is a change - Ve" | "rsion for branch alternate" | ")"
# Target: same three runs, all recolored to FF0000, with the middle run's
# text changed to "rsion for main branch".
#
# Doing the text replacement first (while all three runs share the same
# color) would make this engine's text-mutation path merge the
# now-identically-formatted neighboring runs into one. To keep the three
# runs distinct we first give the middle run a temporarily different
# color (so it cannot merge with its neighbors), perform the text edit,
# and only then unify all three runs' colors to the final FF0000 value
# (a pure formatting change, which does not trigger a run merge).

$middle = $d.Content.Find
$midRange = $d.Range(0, 0)
$midRange.Find.Execute("rsion for branch alternate")
$midRange.Font.Color = 255

$midRange2 = $d.Range(0, 0)
$midRange2.Find.Execute("rsion for branch alternate")
$midRange2.Text = "rsion for main branch"

$run1 = $d.Range(0, 0)
$run1.Find.Execute("(This is a change – Ve")
$run1.Font.Color = 16711680

$run3 = $d.Range(0, 0)
$run3.Find.Execute(")")
$run3.Font.Color = 16711680

$run2 = $d.Range(0, 0)
$run2.Find.Execute("rsion for main branch")
$run2.Font.Color = 16711680
